$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 is the last data row: 04-31 | 삼각김밥 | P9194 | 1200원 | 708
# -> becomes:                05-01 | 휴지     | P8130 | 2700   | 91
$ws.Range("A6").Value = "05-01"
$ws.Range("B6").Value = "휴지"
$ws.Range("C6").Value = "P8130"

# D6/E6 become plain digit strings ("2700", "91"). Excel would normally
# auto-convert those to numbers on entry, but the column stores values as
# text (like "1200원", "1200", "100" elsewhere), so force text formatting
# first to keep them stored as shared-string text rather than numbers.
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "2700"
$ws.Range("E6").Value = "91"
